# :construction: @ExcelStyle 속성 추가
#
# Adds the (default-valued, but now explicit) @ExcelStyle alignment/quote-prefix
# attributes to the two header/label cell styles, and refreshes the sample
# numeric data in columns C:E for rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style tweaks -----------------------------------------------------
# Style used by the header row (A1:E1) -> fontId=1/fillId=4/borderId=8
$headerRange = $ws.Range("A1:E1")
$headerRange.ShrinkToFit = $false
$headerRange.WrapText = $false
$headerRange.Orientation = 0
$headerRange.PrefixCharacter = $false

# Style used by the row-label column (A2:A5) -> fontId=2/fillId=8/borderId=12
$labelRange = $ws.Range("A2:A5")
$labelRange.ShrinkToFit = $false
$labelRange.WrapText = $false
$labelRange.Orientation = 0
$labelRange.PrefixCharacter = $false

# --- Data refresh -------------------------------------------------------
$ws.Range("C2").Value = 760919.0
$ws.Range("D2").Value = 560638.3125
$ws.Range("E2").Value = 35489.0

$ws.Range("C3").Value = 443485.0
$ws.Range("D3").Value = 586050.6875
$ws.Range("E3").Value = 39459.0

$ws.Range("C4").Value = 932822.0
$ws.Range("D4").Value = 934342.625
$ws.Range("E4").Value = 233926.0

$ws.Range("C5").Value = 87801.0
$ws.Range("D5").Value = 566185.0
$ws.Range("E5").Value = 73175.0
